# UI update: modify leaderboard and user manager
#
# - Adds a 3rd ("leaderboard" status/score) column next to the existing
#   username/password-hash table.
# - Adds new user-manager rows (test1, 123, username1, qwe) below the
#   existing admin row.
# - Widens the password-hash column so the long MD5 hashes are readable.
# - Leaves the selection on the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing rows: keep username/password data, add the new 3rd column ---
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "5F4DCC3B5AA765D61D8327DEB882CF99"
$ws.Range("C1").Value = 0

$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "21232F297A57A5A743894A0E4A801FC3"
$ws.Range("C2").Value = 0

# --- user manager: newly added accounts ---
$ws.Range("A3").Value = "test1"
$ws.Range("B3").Value = "d8578edf8458ce06fbc5bb76a58c5ca4"
$ws.Range("C3").Value = 0

# force "123" to be stored as text (a login name), not a number
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "123"
$ws.Range("B4").Value = "202cb962ac59075b964b07152d234b70"
$ws.Range("C4").Value = 0

$ws.Range("A5").Value = "username1"
$ws.Range("B5").Value = "202cb962ac59075b964b07152d234b70"
$ws.Range("C5").Value = 0

$ws.Range("A6").Value = "qwe"
$ws.Range("B6").Value = "202cb962ac59075b964b07152d234b70"
$ws.Range("C6").Value = 0

# --- widen the password-hash column so the hashes are fully visible ---
$ws.Columns("B").ColumnWidth = 43.33203125

# --- move the selection down onto the newly added data ---
[void]$ws.Range("A4:D4").Select()

Write-Output "applied leaderboard/user-manager update"
